$wb = $excel.ActiveWorkbook

$wsMoved   = $wb.Worksheets.Item(1)   # 移到外層的檔名
$wsCond    = $wb.Worksheets.Item(2)   # 有條件使用的檔名
$wsRepeat  = $wb.Worksheets.Item(3)   # 可以重複分配的角度
$wsForbid  = $wb.Worksheets.Item(4)   # 角度禁止規則
$wsKeyword = $wb.Worksheets.Item(5)   # 商品分類及關鍵字條件

# --- Sheet1: 移到外層的檔名 ---
$wsMoved.Range("A1").Value = "檔名含有以下關鍵字者不編"

# --- Sheet2: 有條件使用的檔名 --- add two new rows of data
$wsCond.Range("A5").Value = "_FC_"
$wsCond.Range("B5").Value = "_FR_Torso"
$wsCond.Range("A6").Value = "_F_Model_"
$wsCond.Range("B6").Value = "_FR_Model"

# --- Sheet4: 角度禁止規則 ---
$wsForbid.Range("C1").Value = "禁止邏輯 (等於代表名稱完全相同 ; 包含代表名稱該字串包含即可)"

# --- Sheet5: 商品分類及關鍵字條件 ---
$wsKeyword.Range("B1").Value = "必須含有以下關鍵字才能判定為左欄的商品分類"
$wsKeyword.Range("C1").Value = "是否需要關鍵字全部滿足 (False代表滿足其一即可)"
$wsKeyword.Range("A3").Value = "套裝"
$wsKeyword.Range("B3").Value = "_Btp_,_Fbp_,_Bbp_"
$wsKeyword.Range("C3").Value = $false
